$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new config rows
$ws.Range("A12").Value = "remove_column"
$ws.Range("B12").Value = "Column-10,Last 7 Days,Column-0,Circulating Supply,Volume(24h)"

$ws.Range("A13").Value = "add_column"
$ws.Range("B13").Value = "Old Position,PriceChange,BestBuy,BestSell"

$ws.Range("A14").Value = "max_coin"
$ws.Range("B14").Value = 5

# Widen column B to fit the new content
$ws.Columns.Item(2).ColumnWidth = 56

# Select B14 as the active cell, matching the final saved state
$ws.Range("B14").Select()
